$wb = $excel.ActiveWorkbook

# --- Sheet 1: "3ASY01_RNASeq" (the main annotation table) ---
$ws = $wb.Worksheets.Item("3ASY01_RNASeq")

# Rename the table column header "Parameter [Base-calling Software version]"
# -> "Parameter [Base-calling Software Version]" (capital V). This also
# renames the underlying ListObject column / table definition.
$ws.Range("AG1").Value2 = "Parameter [Base-calling Software Version]"

# Row 2 (Sample row "RNA-Seq strategy")
$ws.Range("C2").Value2  = "user-specific"
$ws.Range("D2").Value2  = "user-specific"
$ws.Range("F2").Value2  = "user-specific"
$ws.Range("G2").Value2  = "user-specific"
$ws.Range("H2").Value2  = "http://purl.obolibrary.org/obo/NFDI4PSO_1000009"
$ws.Range("I2").Value2  = "NFDI4PSO"
$ws.Range("J2").Value2  = "paired-end"
$ws.Range("L2").Value2  = "user-specific"
$ws.Range("M2").Value2  = "user-specific"
$ws.Range("R2").Value2  = "user-specific"
$ws.Range("S2").Value2  = "user-specific"
$ws.Range("U2").Value2  = "http://purl.obolibrary.org/obo/NFDI4PSO_1000087"
$ws.Range("V2").Value2  = "microgram"
$ws.Range("W2").Value2  = "UO"
$ws.Range("AB2").Value2 = "user-specific"
$ws.Range("AC2").Value2 = "user-specific"
$ws.Range("AE2").Value2 = "user-specific"
$ws.Range("AF2").Value2 = "user-specific"
$ws.Range("AH2").Value2 = "user-specific"
$ws.Range("AI2").Value2 = "user-specific"
$ws.Range("AN2").Value2 = "user-specific"
$ws.Range("AO2").Value2 = "user-specific"
$ws.Range("AQ2").Value2 = "user-specific"
$ws.Range("AR2").Value2 = "user-specific"

# Row 3 (Sample row "ChIP-Seq strategy")
$ws.Range("C3").Value2  = "user-specific"
$ws.Range("D3").Value2  = "user-specific"
$ws.Range("F3").Value2  = "user-specific"
$ws.Range("G3").Value2  = "user-specific"
$ws.Range("I3").Value2  = "user-specific"
$ws.Range("J3").Value2  = "user-specific"
$ws.Range("R3").Value2  = "user-specific"
$ws.Range("S3").Value2  = "user-specific"
$ws.Range("U3").Value2  = "http://purl.obolibrary.org/obo/NFDI4PSO_1000087"
$ws.Range("V3").Value2  = "microgram"
$ws.Range("W3").Value2  = "UO"
$ws.Range("AN3").Value2 = "user-specific"
$ws.Range("AO3").Value2 = "user-specific"
$ws.Range("AQ3").Value2 = "user-specific"
$ws.Range("AR3").Value2 = "user-specific"

# Row 4 (Sample row "ssRNA-seq")
$ws.Range("C4").Value2  = "NFDI4PSO"
$ws.Range("D4").Value2  = "1.1.5"
$ws.Range("F4").Value2  = "user-specific"
$ws.Range("G4").Value2  = "user-specific"
$ws.Range("R4").Value2  = "user-specific"
$ws.Range("S4").Value2  = "user-specific"
$ws.Range("U4").Value2  = "http://purl.obolibrary.org/obo/NFDI4PSO_1000087"
$ws.Range("V4").Value2  = "microgram"
$ws.Range("W4").Value2  = "UO"
$ws.Range("AN4").Value2 = "user-specific"
$ws.Range("AO4").Value2 = "user-specific"

# Row 5 (Sample row "rRNA depletion")
$ws.Range("F5").Value2  = "NFDI4PSO"
$ws.Range("G5").Value2  = "http://purl.obolibrary.org/obo/NFDI4PSO_0000082"
$ws.Range("U5").Value2  = "http://purl.obolibrary.org/obo/NFDI4PSO_1000087"
$ws.Range("V5").Value2  = "microgram"
$ws.Range("W5").Value2  = "UO"

# --- Sheet 2: "SwateTemplateMetadata" ---
$ws2 = $wb.Worksheets.Item("SwateTemplateMetadata")
# Template version bump 1.1.5 -> 1.1.6
$ws2.Range("B3").Value2 = "1.1.6"
